$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style source cells from existing data rows (row 2): A col style=1 (bold/border/center), E col style=2 (datetime numfmt)
$styleColA = $ws.Range("A2")
$styleColE = $ws.Range("E2")

# Row 73
$ws.Range("A73").Value = 72
$styleColA.Copy()
$ws.Range("A73").PasteSpecial(-4122)
$ws.Range("B73").Value = "azerbaijan"
$ws.Range("C73").Value = "premier-league"
$ws.Range("D73").Value = "2023-2024"
$ws.Range("E73").Value = 45263.5
$styleColE.Copy()
$ws.Range("E73").PasteSpecial(-4122)
$ws.Range("F73").Value = "Kapaz"
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = "Sabah Baku"
$ws.Range("I73").Value = 2
$ws.Range("J73").Value = 6.39
$ws.Range("K73").Value = "02/12/2023 00:12"
$ws.Range("L73").Value = 6.14
$ws.Range("M73").Value = "03/12/2023 11:58"
$ws.Range("N73").Value = 3.92
$ws.Range("O73").Value = "02/12/2023 00:12"
$ws.Range("P73").Value = 3.69
$ws.Range("Q73").Value = "03/12/2023 11:58"
$ws.Range("R73").Value = 1.44
$ws.Range("S73").Value = "02/12/2023 00:12"
$ws.Range("T73").Value = 1.57
$ws.Range("U73").Value = "03/12/2023 11:58"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/kapaz-sabah-baku/SYfiuokG/"

# Row 74
$ws.Range("A74").Value = 73
$styleColA.Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("B74").Value = "azerbaijan"
$ws.Range("C74").Value = "premier-league"
$ws.Range("D74").Value = "2023-2024"
$ws.Range("E74").Value = 45263.58333333334
$styleColE.Copy()
$ws.Range("E74").PasteSpecial(-4122)
$ws.Range("F74").Value = "Neftci Baku"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Zira"
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2.08
$ws.Range("K74").Value = "02/12/2023 02:13"
$ws.Range("L74").Value = 2.37
$ws.Range("M74").Value = "03/12/2023 13:27"
$ws.Range("N74").Value = 2.88
$ws.Range("O74").Value = "02/12/2023 02:13"
$ws.Range("P74").Value = 2.87
$ws.Range("Q74").Value = "03/12/2023 13:27"
$ws.Range("R74").Value = 3.6
$ws.Range("S74").Value = "02/12/2023 02:13"
$ws.Range("T74").Value = 3.34
$ws.Range("U74").Value = "03/12/2023 13:27"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/neftci-baku-zira-fk/bcBTLdle/"

# Row 75
$ws.Range("A75").Value = 74
$styleColA.Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Range("B75").Value = "azerbaijan"
$ws.Range("C75").Value = "premier-league"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value = 45264.6875
$styleColE.Copy()
$ws.Range("E75").PasteSpecial(-4122)
$ws.Range("F75").Value = "Qarabag"
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = "Gabala"
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1.25
$ws.Range("K75").Value = "03/12/2023 13:43"
$ws.Range("L75").Value = 1.26
$ws.Range("M75").Value = "04/12/2023 16:22"
$ws.Range("N75").Value = 5.26
$ws.Range("O75").Value = "03/12/2023 13:43"
$ws.Range("P75").Value = 5.47
$ws.Range("Q75").Value = "04/12/2023 16:27"
$ws.Range("R75").Value = 8.35
$ws.Range("S75").Value = "03/12/2023 13:43"
$ws.Range("T75").Value = 11.16
$ws.Range("U75").Value = "04/12/2023 16:27"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/qarabag-agdam-gabala/4bIKNIJr/"

# Row 76
$ws.Range("A76").Value = 75
$styleColA.Copy()
$ws.Range("A76").PasteSpecial(-4122)
$ws.Range("B76").Value = "azerbaijan"
$ws.Range("C76").Value = "premier-league"
$ws.Range("D76").Value = "2023-2024"
$ws.Range("E76").Value = 45268.6875
$styleColE.Copy()
$ws.Range("E76").PasteSpecial(-4122)
$ws.Range("F76").Value = "Sumqayit"
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = "Qarabag"
$ws.Range("I76").Value = 6
$ws.Range("J76").Value = 8.06
$ws.Range("K76").Value = "07/12/2023 04:43"
$ws.Range("L76").Value = 8.630000000000001
$ws.Range("M76").Value = "08/12/2023 16:28"
$ws.Range("N76").Value = 4.93
$ws.Range("O76").Value = "07/12/2023 04:43"
$ws.Range("P76").Value = 4.81
$ws.Range("Q76").Value = "08/12/2023 16:28"
$ws.Range("R76").Value = 1.29
$ws.Range("S76").Value = "07/12/2023 04:43"
$ws.Range("T76").Value = 1.34
$ws.Range("U76").Value = "08/12/2023 16:28"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sumqayit-fk-qarabag-agdam/8ShKrKe8/"

# Row 77
$ws.Range("A77").Value = 76
$styleColA.Copy()
$ws.Range("A77").PasteSpecial(-4122)
$ws.Range("B77").Value = "azerbaijan"
$ws.Range("C77").Value = "premier-league"
$ws.Range("D77").Value = "2023-2024"
$ws.Range("E77").Value = 45269.5
$styleColE.Copy()
$ws.Range("E77").PasteSpecial(-4122)
$ws.Range("F77").Value = "Gabala"
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = "Kapaz"
$ws.Range("I77").Value = 3
$ws.Range("J77").Value = 1.79
$ws.Range("K77").Value = "08/12/2023 00:12"
$ws.Range("L77").Value = 1.82
$ws.Range("M77").Value = "09/12/2023 11:50"
$ws.Range("N77").Value = 3.23
$ws.Range("O77").Value = "08/12/2023 00:12"
$ws.Range("P77").Value = 3.32
$ws.Range("Q77").Value = "09/12/2023 11:55"
$ws.Range("R77").Value = 4.19
$ws.Range("S77").Value = "08/12/2023 00:12"
$ws.Range("T77").Value = 4.52
$ws.Range("U77").Value = "09/12/2023 11:18"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/gabala-kapaz/nJiOsvBE/"

# Row 78
$ws.Range("A78").Value = 77
$styleColA.Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("B78").Value = "azerbaijan"
$ws.Range("C78").Value = "premier-league"
$ws.Range("D78").Value = "2023-2024"
$ws.Range("E78").Value = 45269.58333333334
$styleColE.Copy()
$ws.Range("E78").PasteSpecial(-4122)
$ws.Range("F78").Value = "Zira"
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = "Araz"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2.58
$ws.Range("K78").Value = "08/12/2023 02:12"
$ws.Range("L78").Value = 2.22
$ws.Range("M78").Value = "09/12/2023 13:44"
$ws.Range("N78").Value = 2.65
$ws.Range("O78").Value = "08/12/2023 02:12"
$ws.Range("P78").Value = 2.95
$ws.Range("Q78").Value = "09/12/2023 13:44"
$ws.Range("R78").Value = 2.95
$ws.Range("S78").Value = "08/12/2023 02:12"
$ws.Range("T78").Value = 3.54
$ws.Range("U78").Value = "09/12/2023 13:44"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/zira-fk-araz-pfk/hKAXKG31/"

# Row 79
$ws.Range("A79").Value = 78
$styleColA.Copy()
$ws.Range("A79").PasteSpecial(-4122)
$ws.Range("B79").Value = "azerbaijan"
$ws.Range("C79").Value = "premier-league"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45270.5
$styleColE.Copy()
$ws.Range("E79").PasteSpecial(-4122)
$ws.Range("F79").Value = "Sabah Baku"
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = "Sabail"
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 1.4
$ws.Range("K79").Value = "09/12/2023 00:13"
$ws.Range("L79").Value = 1.42
$ws.Range("M79").Value = "10/12/2023 11:16"
$ws.Range("N79").Value = 4.2
$ws.Range("O79").Value = "09/12/2023 00:13"
$ws.Range("P79").Value = 4.39
$ws.Range("Q79").Value = "10/12/2023 11:52"
$ws.Range("R79").Value = 6.64
$ws.Range("S79").Value = "09/12/2023 00:13"
$ws.Range("T79").Value = 7.22
$ws.Range("U79").Value = "10/12/2023 11:52"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sabah-baku-sabail/ttiStbQK/"

# Row 80
$ws.Range("A80").Value = 79
$styleColA.Copy()
$ws.Range("A80").PasteSpecial(-4122)
$ws.Range("B80").Value = "azerbaijan"
$ws.Range("C80").Value = "premier-league"
$ws.Range("D80").Value = "2023-2024"
$ws.Range("E80").Value = 45270.625
$styleColE.Copy()
$ws.Range("E80").PasteSpecial(-4122)
$ws.Range("F80").Value = "Neftci Baku"
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = "Turan"
$ws.Range("I80").Value = 2
$ws.Range("J80").Value = 1.69
$ws.Range("K80").Value = "09/12/2023 03:12"
$ws.Range("L80").Value = 1.72
$ws.Range("M80").Value = "10/12/2023 14:56"
$ws.Range("N80").Value = 3.48
$ws.Range("O80").Value = "09/12/2023 03:12"
$ws.Range("P80").Value = 3.55
$ws.Range("Q80").Value = "10/12/2023 14:56"
$ws.Range("R80").Value = 4.35
$ws.Range("S80").Value = "09/12/2023 03:12"
$ws.Range("T80").Value = 4.75
$ws.Range("U80").Value = "10/12/2023 14:56"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/neftci-baku-turan/WA9yKzJ7/"

# Row 81
$ws.Range("A81").Value = 80
$styleColA.Copy()
$ws.Range("A81").PasteSpecial(-4122)
$ws.Range("B81").Value = "azerbaijan"
$ws.Range("C81").Value = "premier-league"
$ws.Range("D81").Value = "2023-2024"
$ws.Range("E81").Value = 45274.52083333334
$styleColE.Copy()
$ws.Range("E81").PasteSpecial(-4122)
$ws.Range("F81").Value = "Kapaz"
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = "Sumqayit"
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 2.63
$ws.Range("K81").Value = "13/12/2023 00:42"
$ws.Range("L81").Value = 2.53
$ws.Range("M81").Value = "14/12/2023 12:21"
$ws.Range("N81").Value = 2.78
$ws.Range("O81").Value = "13/12/2023 00:42"
$ws.Range("P81").Value = 3.34
$ws.Range("Q81").Value = "14/12/2023 12:21"
$ws.Range("R81").Value = 2.73
$ws.Range("S81").Value = "13/12/2023 00:42"
$ws.Range("T81").Value = 2.65
$ws.Range("U81").Value = "14/12/2023 12:21"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/kapaz-sumqayit-fk/vHW2zd3l/"

# Row 82
$ws.Range("A82").Value = 81
$styleColA.Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("B82").Value = "azerbaijan"
$ws.Range("C82").Value = "premier-league"
$ws.Range("D82").Value = "2023-2024"
$ws.Range("E82").Value = 45275.41666666666
$styleColE.Copy()
$ws.Range("E82").PasteSpecial(-4122)
$ws.Range("F82").Value = "Sabail"
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = "Gabala"
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2.8
$ws.Range("K82").Value = "13/12/2023 22:12"
$ws.Range("L82").Value = 3.17
$ws.Range("M82").Value = "15/12/2023 09:59"
$ws.Range("N82").Value = 2.96
$ws.Range("O82").Value = "13/12/2023 22:12"
$ws.Range("P82").Value = 3.05
$ws.Range("Q82").Value = "15/12/2023 09:59"
$ws.Range("R82").Value = 2.43
$ws.Range("S82").Value = "13/12/2023 22:12"
$ws.Range("T82").Value = 2.33
$ws.Range("U82").Value = "15/12/2023 09:59"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sabail-gabala/0QXbyxlr/"

# Row 83
$ws.Range("A83").Value = 82
$styleColA.Copy()
$ws.Range("A83").PasteSpecial(-4122)
$ws.Range("B83").Value = "azerbaijan"
$ws.Range("C83").Value = "premier-league"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45275.5
$styleColE.Copy()
$ws.Range("E83").PasteSpecial(-4122)
$ws.Range("F83").Value = "Araz"
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = "Sabah Baku"
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2.79
$ws.Range("K83").Value = "14/12/2023 00:12"
$ws.Range("L83").Value = 3.48
$ws.Range("M83").Value = "15/12/2023 11:57"
$ws.Range("N83").Value = 3.17
$ws.Range("O83").Value = "14/12/2023 00:12"
$ws.Range("P83").Value = 3.54
$ws.Range("Q83").Value = "15/12/2023 11:58"
$ws.Range("R83").Value = 2.32
$ws.Range("S83").Value = "14/12/2023 00:12"
$ws.Range("T83").Value = 2
$ws.Range("U83").Value = "15/12/2023 11:57"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/araz-pfk-sabah-baku/fqeWuIuR/"

# Row 84
$ws.Range("A84").Value = 83
$styleColA.Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("B84").Value = "azerbaijan"
$ws.Range("C84").Value = "premier-league"
$ws.Range("D84").Value = "2023-2024"
$ws.Range("E84").Value = 45276.625
$styleColE.Copy()
$ws.Range("E84").PasteSpecial(-4122)
$ws.Range("F84").Value = "Turan"
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = "Zira"
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3.12
$ws.Range("K84").Value = "15/12/2023 03:12"
$ws.Range("L84").Value = 2.75
$ws.Range("M84").Value = "16/12/2023 14:09"
$ws.Range("N84").Value = 2.75
$ws.Range("O84").Value = "15/12/2023 03:12"
$ws.Range("P84").Value = 2.97
$ws.Range("Q84").Value = "16/12/2023 13:02"
$ws.Range("R84").Value = 2.37
$ws.Range("S84").Value = "15/12/2023 03:12"
$ws.Range("T84").Value = 2.7
$ws.Range("U84").Value = "16/12/2023 14:09"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/turan-zira-fk/QXUAYZX0/"

# Row 85
$ws.Range("A85").Value = 84
$styleColA.Copy()
$ws.Range("A85").PasteSpecial(-4122)
$ws.Range("B85").Value = "azerbaijan"
$ws.Range("C85").Value = "premier-league"
$ws.Range("D85").Value = "2023-2024"
$ws.Range("E85").Value = 45278.70833333334
$styleColE.Copy()
$ws.Range("E85").PasteSpecial(-4122)
$ws.Range("F85").Value = "Qarabag"
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = "Neftci Baku"
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1.44
$ws.Range("K85").Value = "17/12/2023 05:13"
$ws.Range("L85").Value = 1.25
$ws.Range("M85").Value = "18/12/2023 16:20"
$ws.Range("N85").Value = 4.16
$ws.Range("O85").Value = "17/12/2023 05:13"
$ws.Range("P85").Value = 5.67
$ws.Range("Q85").Value = "18/12/2023 16:59"
$ws.Range("R85").Value = 5.91
$ws.Range("S85").Value = "17/12/2023 05:13"
$ws.Range("T85").Value = 11.27
$ws.Range("U85").Value = "18/12/2023 16:59"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/qarabag-agdam-neftci-baku/nsW6ZFIf/"

